# ------------------------------------------------------------------
# Sh_report.xlsx edit: split single "Sheet1" workbook into three tabs
# (Swing / Shortlist / Watchlist), add a Watchlist stock table with an
# AutoFilter on the "Type" column, and add an "Indicators" header to
# the Swing sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet to "Swing" -----------------------
$swing = $wb.Worksheets.Item(1)
$swing.Name = "Swing"

# --- 2. Add the two new sheets in the right order / sheetId scheme -
# Real Excel assigns sheetId in *creation* order but keeps tabs in
# the order they were positioned. Creating Watchlist first (so it
# gets the lower sheetId) then Shortlist (higher sheetId), both
# inserted right after Swing, reproduces:
#   Swing(sheetId1) / Shortlist(sheetId3) / Watchlist(sheetId2)
$tmp = $wb.Worksheets.Add($null, $swing)
$tmp.Name = "Watchlist"

$tmp2 = $wb.Worksheets.Add($null, $swing)
$tmp2.Name = "Shortlist"

# NOTE: sheet object handles returned by Add() track *position*, not
# identity - once a later Add() shifts a sheet's index, an old handle
# silently starts resolving to whatever sheet now sits at that index.
# Re-fetch fresh handles by name now that the sheet collection is
# stable, and use only these from here on.
$swing     = $wb.Worksheets.Item("Swing")
$shortlist = $wb.Worksheets.Item("Shortlist")
$watchlist = $wb.Worksheets.Item("Watchlist")

# --- 3. Swing sheet: add new "Indicators" header in column O -------
$swing.Range("O1").Value = "Indicators"

# --- 4. Watchlist sheet: header row -------------------------------
$watchlist.Range("A2").Value = "Date"
$watchlist.Range("B2").Value = "Shares"
$watchlist.Range("C2").Value = "LTP"
$watchlist.Range("D2").Value = "Target"
$watchlist.Range("E2").Value = "percentage"
$watchlist.Range("F2").Value = "StopLoss"
$watchlist.Range("G2").Value = "percentage"
$watchlist.Range("H2").Value = "Type"
$watchlist.Range("I2").Value = "MACD"
$watchlist.Range("J2").Value = "SRSI"
$watchlist.Range("K2").Value = "S/R"
$watchlist.Range("L2").Value = "Bband"
$watchlist.Range("M2").Value = "RSI"
$watchlist.Range("N2").Value = "ADX"
$watchlist.Range("O2").Value = "HA"

# --- 5. Watchlist sheet: data rows ---------------------------------
# Row 3 - BAJAJ-AUTO
$watchlist.Range("A3").Value = "26/05/2019"
$watchlist.Range("B3").Value = "BAJAJ-AUTO"
$watchlist.Range("C3").Value = 3079.35
$watchlist.Range("D3").Value = 3171
$watchlist.Range("E3").Value = 3
$watchlist.Range("F3").Value = 2986
$watchlist.Range("G3").Value = 3
$watchlist.Range("H3").Value = "S1"
$watchlist.Range("I3").Value = "yes"
$watchlist.Range("J3").Value = "yes"
$watchlist.Range("L3").Value = "yes"
$watchlist.Range("O3").Value = "yes"

# Row 4 - JSWSTEEL
$watchlist.Range("A4").Value = "26/05/2019"
$watchlist.Range("B4").Value = "JSWSTEEL"
$watchlist.Range("H4").Value = "S5"
$watchlist.Range("I4").Value = "yes"
$watchlist.Range("L4").Value = "yes"
$watchlist.Range("O4").Value = "yes"

# Row 5 - CIPLA
$watchlist.Range("A5").Value = "26/05/2019"
$watchlist.Range("B5").Value = "CIPLA"
$watchlist.Range("C5").Value = 570.3
$watchlist.Range("H5").Value = "S2"
$watchlist.Range("I5").Value = "yes"
$watchlist.Range("L5").Value = "yes"
$watchlist.Range("O5").Value = "yes"

# Row 6 - COALINDIA
$watchlist.Range("A6").Value = "26/05/2019"
$watchlist.Range("B6").Value = "COALINDIA"
$watchlist.Range("H6").Value = "S5"
$watchlist.Range("I6").Value = "yes"
$watchlist.Range("L6").Value = "yes"
$watchlist.Range("O6").Value = "yes"

# Row 7 - GRASIM
$watchlist.Range("A7").Value = "26/05/2019"
$watchlist.Range("B7").Value = "GRASIM"
$watchlist.Range("C7").Value = 907.1
$watchlist.Range("H7").Value = "S3"
$watchlist.Range("I7").Value = "yes"
$watchlist.Range("L7").Value = "yes"
$watchlist.Range("O7").Value = "yes"

# Row 8 - M&M
$watchlist.Range("A8").Value = "26/05/2019"
$watchlist.Range("B8").Value = "M&M"
$watchlist.Range("C8").Value = 665.7
$watchlist.Range("H8").Value = "S3"
$watchlist.Range("I8").Value = "yes"
$watchlist.Range("L8").Value = "yes"
$watchlist.Range("O8").Value = "yes"

# Row 9 - NTPC
$watchlist.Range("A9").Value = "26/05/2019"
$watchlist.Range("B9").Value = "NTPC"
$watchlist.Range("H9").Value = "S5"
$watchlist.Range("I9").Value = "yes"
$watchlist.Range("L9").Value = "yes"
$watchlist.Range("O9").Value = "yes"

# Row 10 - YESBANK
$watchlist.Range("A10").Value = "26/05/2019"
$watchlist.Range("B10").Value = "YESBANK"
$watchlist.Range("H10").Value = "S4"
$watchlist.Range("I10").Value = "yes"
$watchlist.Range("L10").Value = "yes"
$watchlist.Range("O10").Value = "yes"

# Row 11 - INFRATEL
$watchlist.Range("A11").Value = "26/05/2019"
$watchlist.Range("B11").Value = "INFRATEL"
$watchlist.Range("C11").Value = 278.6
$watchlist.Range("H11").Value = "S2"
$watchlist.Range("I11").Value = "yes"
$watchlist.Range("L11").Value = "yes"
$watchlist.Range("O11").Value = "yes"

# Row 12 - POWERGRID
$watchlist.Range("A12").Value = "26/05/2019"
$watchlist.Range("B12").Value = "POWERGRID"
$watchlist.Range("H12").Value = "S5"
$watchlist.Range("I12").Value = "yes"
$watchlist.Range("L12").Value = "yes"
$watchlist.Range("O12").Value = "yes"

# --- 6. Column widths on Watchlist (bestFit look of the original) --
$watchlist.Columns.Item(1).ColumnWidth = 10.7109375
$watchlist.Columns.Item(2).ColumnWidth = 11.85546875
$watchlist.Columns.Item(3).ColumnWidth = 16.85546875
$watchlist.Columns.Item(4).ColumnWidth = 8.85546875
$watchlist.Columns.Item(5).ColumnWidth = 13.28515625
$watchlist.Columns.Item(6).ColumnWidth = 11
$watchlist.Columns.Item(7).ColumnWidth = 11
$watchlist.Columns.Item(9).ColumnWidth = 15.140625

# --- 7. AutoFilter on Type column (H, 8th col of A:K) ---------------
$filterRange = $watchlist.Range("A2:K12")
$filterRange.AutoFilter(8, @("S1", "S2", "S3"), 7)

# --- 8. Selections / active sheet -----------------------------------
$shortlist.Range("I15").Select()
$swing.Range("L7").Select()

$watchlist.Range("E15").Select()
$watchlist.Select()

$wb.Windows.Item(1).WindowState = -4143
